$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The match stats in row 2 and row 3 were swapped (runs, balls, sixes);
# "fours" (column E) stayed 0 in both rows so it is unaffected.
# Force text number format on each changed cell so the value keeps being
# stored as a text string (matching the original string-typed cells)
# rather than being auto-converted to a number by Excel.

$ws.Range("C2").NumberFormat = "@"
$ws.Range("C2").Value = "10"

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "8"

$ws.Range("F2").NumberFormat = "@"
$ws.Range("F2").Value = "1"

$ws.Range("C3").NumberFormat = "@"
$ws.Range("C3").Value = "2"

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3"

$ws.Range("F3").NumberFormat = "@"
$ws.Range("F3").Value = "0"
